# Swap the species-observation data between row 5 and row 6.
# Columns that differ between the two rows (A, B, E, F, G, H, I, Q, R, AC)
# need to have their values exchanged; all other columns stay as-is.
# Note: use Value2 (not Value) when reading into a variable, since Value
# round-trips poorly through an intermediate PowerShell variable here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 5
$row2 = 6

# Columns whose values differ between the two rows and must be swapped.
$cols = @("A", "B", "E", "F", "G", "H", "I", "Q", "R", "AC")

foreach ($col in $cols) {
    $rng1 = $ws.Range("$col$row1")
    $rng2 = $ws.Range("$col$row2")

    $val1 = $rng1.Value2
    $val2 = $rng2.Value2

    $rng1.Value2 = $val2
    $rng2.Value2 = $val1
}
